$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "1.986.58"
$ws.Range("E3").Value = "  -3.63%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.46"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -4.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.43"
$ws.Range("E7").Value = "  +7.29%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.62"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.365"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0741"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.952"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.72"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "2.276.14"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.86"
$ws.Range("E17").Value = "  +7.64%  "
$ws.Range("D18").Value = "1.985.74"
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("D19").Value = "35.731.10"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.83"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.21"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.56"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +13.93%  "
$ws.Range("E26").Value = "  -4.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.79"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.33"
$ws.Range("E29").Value = "  -4.32%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.92"
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("E32").Value = "  -6.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0977"
$ws.Range("E33").Value = "  +15.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0602"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  +11.25%  "
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.70"
$ws.Range("E39").Value = "  +12.29%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0940"
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0214"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.54"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "94.11"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.76"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  -3.42%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.18"
$ws.Range("E51").Value = "  +2.63%  "
